# Update workbook to add two new Pick&Place rows (R1, R2) and rename the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (2024-12-13 -> 2024-12-21)
$ws.Name = "PickAndPlace_PCB1_2024-12-21"

# Row 8 - R1
$ws.Range("A8").Value = "R1"
$ws.Range("B8").Value = "0805W8F1002T5E"
$ws.Range("C8").Value = "R0805"
$ws.Range("D8").Value = "62.103mm"
$ws.Range("E8").Value = "21.336mm"
$ws.Range("F8").Value = "62.103mm"
$ws.Range("G8").Value = "21.336mm"
$ws.Range("H8").Value = "61.103mm"
$ws.Range("I8").Value = "21.336mm"
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = "T"
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = "Yes"
$ws.Range("N8").Value = "10kΩ"

# Row 9 - R2
$ws.Range("A9").Value = "R2"
$ws.Range("B9").Value = "0805W8F1002T5E"
$ws.Range("C9").Value = "R0805"
$ws.Range("D9").Value = "61.781mm"
$ws.Range("E9").Value = "12.348mm"
$ws.Range("F9").Value = "61.781mm"
$ws.Range("G9").Value = "12.348mm"
$ws.Range("H9").Value = "60.781mm"
$ws.Range("I9").Value = "12.348mm"
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = "T"
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = "Yes"
$ws.Range("N9").Value = "10kΩ"
